$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 11.7
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 10.4
$ws.Range("B4").Value = 0.8
$ws.Range("C5").Value = 17
